$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "66.306.51"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +3.93%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.489.11"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +2.21%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "595.48"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +4.51%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "169.44"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +7.66%  "
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.488.19"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +2.14%  "
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +1.56%  "
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +0.71%  "
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +4.14%  "
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +1.94%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.096.99"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +2.41%  "
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +0.54%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "27.90"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +2.95%  "
$ws.Range("B16").NumberFormat = "@"
$ws.Range("B16").Value = "WrappedBTC"
$ws.Range("C16").NumberFormat = "@"
$ws.Range("C16").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "66.250.81"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +3.71%  "
$ws.Range("B17").NumberFormat = "@"
$ws.Range("B17").Value = "ShibaInu"
$ws.Range("C17").NumberFormat = "@"
$ws.Range("C17").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000176"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +1.79%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.505.85"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +2.97%  "
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +2.46%  "
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +3.02%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "387.25"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +2.63%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.99"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +2.97%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "72.93"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +2.46%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.999"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -0.10%  "
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +1.70%  "
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +5.15%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.12"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +4.33%  "
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +1.37%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.00"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +0.17%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.38"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +5.12%  "
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +5.33%  "
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +4.31%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "23.44"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +2.66%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "7.41"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +6.26%  "
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +0.05%  "
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +1.13%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "160.85"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +0.07%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.902"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +9.04%  "
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +4.86%  "
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +2.62%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "27.26"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +5.45%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "26.46"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +0.91%  "
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +4.43%  "
$ws.Range("B44").NumberFormat = "@"
$ws.Range("B44").Value = "Filecoin"
$ws.Range("C44").NumberFormat = "@"
$ws.Range("C44").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "4.56"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +3.34%  "
$ws.Range("B45").NumberFormat = "@"
$ws.Range("B45").Value = "Maker"
$ws.Range("C45").NumberFormat = "@"
$ws.Range("C45").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.804.00"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -0.28%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "43.38"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +1.52%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0313"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +2.31%  "
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +6.11%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "349.55"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +8.18%  "
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +5.26%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "32.77"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +7.82%  "
